# "separate fields into children scripts"
# Append a new timesheet entry (row 53) following the same pattern as the
# existing rows: date in col A, hours worked in col B, running total
# formula in col C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prevRow = 52
$newRow  = 53

# A53: date (serial 45410 == 2024-04-28). Copy the previous date cell's
# formatting first so the new cell reuses the existing date style instead
# of Excel minting a brand new number format for it.
$ws.Cells.Item($prevRow, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item($newRow, 1).Value = 45410

# B53: hours worked that day.
$ws.Cells.Item($newRow, 2).Value = 2

# C53: running total, same "previous total + today's hours" formula used
# by every other row in column C.
$ws.Cells.Item($newRow, 3).Formula = "=C52+B53"

# Keep the live selection on the newest row, same as before the edit.
$ws.Range("C53").Activate()
